# Add a new "2022-Q1" sheet (fund-holding detail) right before the "总计"
# (totals) sheet, and prepend a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet named "2022-Q1" immediately before "总计"
# ---------------------------------------------------------------------
$totalSheetRef = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetRef, [System.Type]::Missing)
$newSheet.Name = "2022-Q1"

# NOTE: a worksheet reference captured before Add() tracks a sheet
# *position* rather than the sheet object itself, so after inserting a
# new sheet ahead of "总计" we must re-resolve it by name to get a
# reference that actually points at the totals sheet from now on.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$hdr = $newSheet.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data rows (codes / ratios kept as text so leading zeros survive, just
# like the other quarter sheets in this workbook)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'952004"
$newSheet.Range("C2").Value = "'国泰君安君得明混合"
$newSheet.Range("D2").Value = "'22.12"
$newSheet.Range("E2").Value = "'76.15"
$newSheet.Range("F2").Value = "'4.07"
$newSheet.Range("G2").Value = "'0.9003"
$newSheet.Range("H2").Value = 2

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001899"
$newSheet.Range("C3").Value = "'东海中证社会发展安全产业主题指数"
$newSheet.Range("D3").Value = "'0.21"
$newSheet.Range("E3").Value = "'90.30"
$newSheet.Range("F3").Value = "'2.87"
$newSheet.Range("G3").Value = "'0.0060"
$newSheet.Range("H3").Value = 4

$idxCol = $newSheet.Range("A2:A3")
$idxCol.Font.Bold = $true
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
$idxCol.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing rows down by one (done manually, bottom-up, instead of
#    Rows.Insert, so every cell keeps its original content/type).
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $rn = $r + 1
    $aVal = $totalSheet.Range("A$r").Value2
    $bVal = $totalSheet.Range("B$r").Value2
    $cVal = $totalSheet.Range("C$r").Value2
    $dVal = $totalSheet.Range("D$r").Value2

    $totalSheet.Range("A$rn").Value2 = $aVal + 1
    $totalSheet.Range("B$rn").Value2 = $bVal
    $totalSheet.Range("C$rn").Value2 = $cVal
    $totalSheet.Range("D$rn").Value2 = $dVal
}

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.91

$a2 = $totalSheet.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
